# Moving from AQL 7.x to 8.x.
#
# Append two extra "<---  <message>" annotation blocks (each preceded by a
# run of 4 spaces) right after the existing error annotation run at the end
# of the second paragraph:
#   "    <---Couldn't find the 'self' variable    <---missing feature access or service call"
#
# New runs re-use the same character formatting as the existing
# "<---"/"Expression ..." runs: red text (FF0000), size 16pt (sz=32
# half-points), highlighted with WdColorIndex.wdGray25 (-> w:highlight
# "lightGray"). The "    " (4-space) separator runs stay unformatted, just
# like the existing one already in the document.

$d = $word.ActiveDocument

# The text we need to append goes at the very end of the 2nd paragraph,
# i.e. right before its paragraph mark.
$p2 = $d.Paragraphs(2)
$insPos = $p2.Range.End - 1

function Insert-PlainRun([int]$pos, [string]$text) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    return $pos + $text.Length
}

function Insert-FormattedRun([int]$pos, [string]$text) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $newEnd = $pos + $text.Length
    $fr = $d.Range($pos, $newEnd)
    $fr.Font.Color = 255
    $fr.Font.Size = 16
    $fr.Font.HighlightColorIndex = 16
    return $newEnd
}

# --- First annotation: "    <---Couldn't find the 'self' variable"
$insPos = Insert-PlainRun $insPos "    "
$insPos = Insert-FormattedRun $insPos "<---"
$insPos = Insert-FormattedRun $insPos "Couldn't find the 'self' variable"

# --- Second annotation: "    <---missing feature access or service call"
$insPos = Insert-PlainRun $insPos "    "
$insPos = Insert-FormattedRun $insPos "<---"
$insPos = Insert-FormattedRun $insPos "missing feature access or service call"
